# "used single ddt and added screenshot with extent report"
# Adds two new worksheets (Contact, Registration) with test data,
# mirroring the existing Login/Checkout DDT sheets' layout & styling.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Add the two new worksheets at the end, in order.
# ---------------------------------------------------------------
$afterCheckout = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsContact = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterCheckout)
$wsContact.Name = "Contact"

$afterContact = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsReg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterContact)
$wsReg.Name = "Registration"

# ---------------------------------------------------------------
# 2. Contact sheet data (values + non-hyperlink styling only)
# ---------------------------------------------------------------
$wsContact.Range("A1").Value = "Errol"
$wsContact.Range("B1").Value = "demo@email.com"
$wsContact.Range("C1").Value = "Testing if this field is working without"

$wsContact.Range("A2").Font.Color = 0
$wsContact.Range("A2").Value = "Tester"
$wsContact.Range("B2").Value = "demo1@example"
$wsContact.Range("C2").Value = "When will the iPad Air be in stock again"

$wsContact.Range("A3").Value = "Te"
$wsContact.Range("B3").Value = "demo1@example.com"
$wsContact.Range("C3").Value = "When will the iPad Air be in stock again"

$wsContact.Range("A4").Value = "Tester"
$wsContact.Range("B4").Value = "demo1@example.com"
$wsContact.Range("C4").Value = "abcd"

# ---------------------------------------------------------------
# 3. Registration sheet data (values + non-hyperlink styling only)
# ---------------------------------------------------------------

# Rows 1-5: default-styled rows
$wsReg.Range("A1").Value = "Errol5"
$wsReg.Range("B1").Value = "C5"
$wsReg.Range("C1").Value = "demo1@example.com"
$wsReg.Range("D1").Value = 12345
$wsReg.Range("E1").Value = "test1234"
$wsReg.Range("F1").Value = "test1234"

$wsReg.Range("A2").Value = "Errol5"
$wsReg.Range("B2").Value = "C5"
$wsReg.Range("C2").Value = "demo2@example.com"
$wsReg.Range("D2").Value = 122345
$wsReg.Range("E2").Value = 123
$wsReg.Range("F2").Value = 123

$wsReg.Range("A3").Value = "Errol5"
$wsReg.Range("C3").Value = "demo3@example.com"
$wsReg.Range("D3").Value = 123455
$wsReg.Range("E3").Value = "test1234"
$wsReg.Range("F3").Value = "test1234"

$wsReg.Range("A4").Value = "Errol5"
$wsReg.Range("B4").Value = "C5"
$wsReg.Range("C4").Value = "demo4@example.com"
$wsReg.Range("D4").Value = 12345
$wsReg.Range("E4").Value = "test1234"
$wsReg.Range("F4").Value = "test1234"

$wsReg.Range("A5").Value = "Errol5"
$wsReg.Range("B5").Value = "C5"
$wsReg.Range("C5").Value = "demo5@example.com"
$wsReg.Range("D5").Value = "abcdefg"
$wsReg.Range("E5").Value = "test1234"
$wsReg.Range("F5").Value = "test1234"

# Rows 6-12: black-font styled rows (fontId 3 -> cellXf applyFont)
$wsReg.Range("A6:F12").Font.Color = 0

$wsReg.Range("A6").Value = "Raunak"
$wsReg.Range("B6").Value = "Naik"
$wsReg.Range("C6").Value = "testing@example.com"
$wsReg.Range("D6").Value = 2123434565
$wsReg.Range("E6").Value = "testing123"
$wsReg.Range("F6").Value = "test4567"

$wsReg.Range("A7").Value = "Raunak"
$wsReg.Range("B7").Value = "Naik"
$wsReg.Range("C7").Value = "demo4@example.com"
$wsReg.Range("D7").Value = 2123434565
$wsReg.Range("E7").Value = "testing123"
$wsReg.Range("F7").Value = "testing123"

$wsReg.Range("B8").Value = "Naik"
$wsReg.Range("C8").Value = "testing1@example.com"
$wsReg.Range("D8").Value = 2123434565
$wsReg.Range("E8").Value = "testing123"
$wsReg.Range("F8").Value = "testing123"

$wsReg.Range("A9").Value = "Raunak"
$wsReg.Range("B9").Value = "Naik"
$wsReg.Range("C9").Value = "testing2@example.com"
$wsReg.Range("D9").NumberFormat = "@"
$wsReg.Range("D9").Value = "1"
$wsReg.Range("E9").Value = "testing123"
$wsReg.Range("F9").Value = "testing123"

$wsReg.Range("A10").Value = "Raunakabcdefghijklmnopqrstuvwxyza"
$wsReg.Range("B10").Value = "Raunak"
$wsReg.Range("C10").Value = "testing3@example.com"
$wsReg.Range("D10").Value = 2123434565
$wsReg.Range("E10").Value = "testing123"
$wsReg.Range("F10").Value = "testing123"

$wsReg.Range("A11").Value = "Raunak"
$wsReg.Range("B11").Value = "Naikabcdefghijklmnopqrstuvwxyzabc"
$wsReg.Range("C11").Value = "testing4@example.com"
$wsReg.Range("D11").Value = 2123434565
$wsReg.Range("E11").Value = "testing123"
$wsReg.Range("F11").Value = "testing123"

$wsReg.Range("A12").Value = "Raunak"
$wsReg.Range("B12").Value = "Naik"
$wsReg.Range("C12").Value = "testing5@example.com"
$wsReg.Range("D12").NumberFormat = "@"
$wsReg.Range("D12").Value = "22222222222222222222222222222222222"
$wsReg.Range("E12").Value = "testing123"
$wsReg.Range("F12").Value = "testing123"

# ---------------------------------------------------------------
# 4. Hyperlinks (added last, after all custom fonts/number formats
#    are registered, so they don't perturb earlier style indices).
# ---------------------------------------------------------------
$wsContact.Hyperlinks.Add($wsContact.Range("B1"), "mailto:demo@email.com")
$wsContact.Hyperlinks.Add($wsContact.Range("B2"), "mailto:demo1@example")
$wsContact.Hyperlinks.Add($wsContact.Range("B3"), "mailto:demo1@example.com")
$wsContact.Hyperlinks.Add($wsContact.Range("B4"), "mailto:demo1@example.com")
$wsContact.Range("B1").Style = "Hyperlink"
$wsContact.Range("B2").Style = "Hyperlink"
$wsContact.Range("B3").Style = "Hyperlink"
$wsContact.Range("B4").Style = "Hyperlink"

$wsReg.Hyperlinks.Add($wsReg.Range("C1"), "mailto:demo1@example.com")
$wsReg.Hyperlinks.Add($wsReg.Range("C2"), "mailto:demo2@example.com")
$wsReg.Hyperlinks.Add($wsReg.Range("C3"), "mailto:demo3@example.com")
$wsReg.Hyperlinks.Add($wsReg.Range("C4"), "mailto:demo4@example.com")
$wsReg.Hyperlinks.Add($wsReg.Range("C5"), "mailto:demo5@example.com")
$wsReg.Hyperlinks.Add($wsReg.Range("C7"), "mailto:demo4@example.com")
$wsReg.Hyperlinks.Add($wsReg.Range("C8"), "mailto:testing1@example.com")
$wsReg.Hyperlinks.Add($wsReg.Range("C9"), "mailto:testing2@example.com")
$wsReg.Hyperlinks.Add($wsReg.Range("C10"), "mailto:testing3@example.com")
$wsReg.Hyperlinks.Add($wsReg.Range("C11"), "mailto:testing4@example.com")
$wsReg.Hyperlinks.Add($wsReg.Range("C12"), "mailto:testing5@example.com")

$wsReg.Range("C1").Style = "Hyperlink"
$wsReg.Range("C2").Style = "Hyperlink"
$wsReg.Range("C3").Style = "Hyperlink"
$wsReg.Range("C4").Style = "Hyperlink"
$wsReg.Range("C5").Style = "Hyperlink"
$wsReg.Range("C7").Style = "Hyperlink"
$wsReg.Range("C8").Style = "Hyperlink"
$wsReg.Range("C9").Style = "Hyperlink"
$wsReg.Range("C10").Style = "Hyperlink"
$wsReg.Range("C11").Style = "Hyperlink"
$wsReg.Range("C12").Style = "Hyperlink"

# ---------------------------------------------------------------
# 5. Selections matching the source workbook's last-saved state.
# ---------------------------------------------------------------
$wsContact.Range("E8").Select()
$wsReg.Range("F14").Select()
